$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "26.11.2020 MC Sales Detils"
#
# Row 2 is repointed from the old "Laxmi Telecom" retailer to the new
# "Mayer Doa" retailer (DSR-0351 / Valugachi / Ruhul Amin / Puthia,
# Rajshahi / ZSO-0023), and a brand-new row 3 is added for the
# "Babu Electronics" retailer (DSR-0247 / Islabari / Md Babu Hosen /
# Natore Sadar, Natore / ZSO-0023).
#
# The individual cell writes below are deliberately ordered (rather
# than going column-by-column) so that newly-introduced shared strings
# land in the workbook's shared-string table in the same order as the
# source edit.
# ------------------------------------------------------------------

# DSR ID that used to live in B2 (DSR-0247) now belongs to the new row 3.
$ws.Cells.Item(3, 2).Value = "DSR-0247"

# Row 2 - new retailer's name / market / owner.
$ws.Cells.Item(2, 3).Value = "Mayer Doa "
$ws.Cells.Item(2, 4).Value = "Valugachi"
$ws.Cells.Item(2, 5).Value = "Ruhul Amin"

# Row 2 - new DSR ID.
$ws.Cells.Item(2, 2).Value = "DSR-0351"

# Row 2 - Thana / District.
$ws.Cells.Item(2, 12).Value = "Puthia"
$ws.Cells.Item(2, 11).Value = "Rajshahi"

# Row 2 - Address.
$ws.Cells.Item(2, 14).Value = "Valugachi, Naopara,Puthia, Rajshahi"

# Row 3 - retailer's name / market / owner.
$ws.Cells.Item(3, 3).Value = "Babu Electronics"
$ws.Cells.Item(3, 4).Value = "Islabari"
$ws.Cells.Item(3, 5).Value = "Md Babu Hosen"

# Row 3 - Thana.
$ws.Cells.Item(3, 12).Value = "Natore Sadar"

# Row 2 - ZSO-ID (also reused by row 3).
$ws.Cells.Item(2, 13).Value = "ZSO-0023"

# Row 3 - Address.
$ws.Cells.Item(3, 14).Value = "Islabari, Natore"

# Row 2 - remaining fields (Contact Person reuses Owner Name, phone numbers).
$ws.Cells.Item(2, 9).Value = "Ruhul Amin"
$ws.Cells.Item(2, 10).Value = 1827156464
$ws.Cells.Item(2, 16).Value = 1827156464
$ws.Cells.Item(2, 20).Value = 1827156464

# Row 3 - remaining fields (Dealer ID / RType / Contact Person / phone /
# District / ZSO-ID / category / geo / transaction method reuse values
# that already exist elsewhere in the sheet).
$ws.Cells.Item(3, 1).Value = "DEL-0179"
$ws.Cells.Item(3, 7).Value = "GO"
$ws.Cells.Item(3, 9).Value = "Md Babu Hosen"
$ws.Cells.Item(3, 10).Value = 1723632345
$ws.Cells.Item(3, 11).Value = "Natore"
$ws.Cells.Item(3, 13).Value = "ZSO-0023"
$ws.Cells.Item(3, 16).Value = 1723632345
$ws.Cells.Item(3, 17).Value = "C"
$ws.Cells.Item(3, 18).Value = "Rural"
$ws.Cells.Item(3, 19).Value = "bKash"
$ws.Cells.Item(3, 20).Value = 1723632345

# Move the live selection to where the author left it.
$ws.Range("S17").Select() | Out-Null
